# Add an "Email" line to the author-info text box on slide 1.
#
# Before:
#   Authors:  Tsz-Yeung Lau
#   Student ID: 11327605
#   Affiliation: Dept. of Computer Science & Information Engineering, Chaoyang University of Technology
#   Date: 07 January 2025
#
# After:
#   Authors:  Tsz-Yeung Lau
#   Student ID: 11327605
#   Email: tylau70242@gmail.com      <-- new paragraph
#   Affiliation: Dept. of Computer Science & Information Engineering, Chaoyang University of Technology
#   Date: 07 January 2025

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the author-info text box robustly (don't hardcode a shape index):
# it's the shape whose text contains "Student ID".
$targetShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $cand = $s.Shapes.Item($i)
    if ($cand.HasTextFrame) {
        if ($cand.TextFrame.HasText) {
            if ($cand.TextFrame.TextRange.Text -like "*Student ID*") {
                $targetShape = $cand
            }
        }
    }
}

$tf = $targetShape.TextFrame

# The "Student ID: 11327605" paragraph is paragraph 2.
$studentIdPara = $tf.TextRange.Paragraphs(2, 1)

# Insert a new paragraph right after it containing "Email" (a paragraph
# break is a carriage return). The new paragraph/run inherits the
# surrounding run formatting (scheme color accent2 / lumMod 75%).
[void]$studentIdPara.InsertAfter([char]13 + "Email")

# Re-acquire the new paragraph (paragraph 3 now) and give its run the
# same sz/latin typeface used by the sibling lines in this box.
$emailPara = $tf.TextRange.Paragraphs(3, 1)
$emailPara.Font.Size = 18
$emailPara.Font.Name = "+mn-lt"

# Append the second run with the email address itself, matching style.
[void]$emailPara.InsertAfter(": tylau70242@gmail.com")
$emailPara2 = $tf.TextRange.Paragraphs(3, 1)
$emailPara2.Font.Size = 18
$emailPara2.Font.Name = "+mn-lt"

# The shape auto-fits its height to the text (spAutoFit); pin it to the
# exact resulting height (EMU -> points, 12700 EMU per point).
$targetShape.Height = 1732141 / 12700
